# Add a new "LastChecked" column before the existing "Matched" column (I),
# shifting Matched/Verdict/Explanation/NeedExplanation one column to the right
# (I->J, J->K, K->L, L->M), then populate the new LastChecked column with the
# date each row was checked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at I; existing I..L shift to J..M.
$ws.Columns.Item(9).Insert()

# Header for the new column.
$ws.Range("I1").Value = "LastChecked"

# Make sure the date values are stored as literal text (not auto-converted
# to Excel date serials) to match "2023-07-09"-style strings.
$ws.Range("I2:I71").NumberFormat = "@"

# Default last-checked date used for almost every row.
for ($r = 2; $r -le 71; $r++) {
    $ws.Cells.Item($r, 9).Value = "2023-07-09"
}

# A couple of rows were checked on different dates.
$ws.Cells.Item(51, 9).Value = "2023-07-01"
$ws.Cells.Item(62, 9).Value = "2023-07-06"
